# "final k vs acc results"
# Adds a third "std" column with a per-k standard deviation formula,
# narrows the spacer column A, repositions/resizes the chart (it was
# dragged/resized in the source edit), and leaves the new data range
# selected, matching the final saved state of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "std" header + STDEVPA formula for the first k (row 2)
$ws.Range("C1").Value = "std"
$ws.Range("C2").Formula = "=STDEVPA(B2:B6)"

# Spacer column A got narrowed
$ws.Columns("A").ColumnWidth = 3.05

# Chart was moved/resized: anchor now from F7 (col 5, offset 213359 EMU;
# row 6, offset 121920 EMU) to S30 (col 18, offset 403858 EMU;
# row 29, offset 182880 EMU). Convert those EMU offsets (12700 EMU/pt)
# onto the worksheet's own column/row pixel positions so the anchor
# lines up exactly regardless of column width quirks.
$co = $ws.ChartObjects().Item(1)
$colFrom = $ws.Columns.Item(6)
$colTo = $ws.Columns.Item(19)
$rowFrom = $ws.Rows.Item(7)
$rowTo = $ws.Rows.Item(30)

$newLeft = $colFrom.Left + (213359 / 12700)
$newTop = $rowFrom.Top + (121920 / 12700)
$newRight = $colTo.Left + (403858 / 12700)
$newBottom = $rowTo.Top + (182880 / 12700)

$co.Left = $newLeft
$co.Top = $newTop
$co.Width = $newRight - $newLeft
$co.Height = $newBottom - $newTop

# Final selection left on the new std column's data range
$ws.Range("C10:D15").Select()
